# pcc_list_by_year.xlsx - "diff in diff starter script"
#
# The "changes" sheet previously had an AutoFilter applied on column C
# ("Labour" only), which both hid the non-matching rows and recorded the
# filter criteria in the workbook. This edit clears that filter (restoring
# all rows to view) and adds a new "change type" analysis column (H) that
# mirrors the "No change" verdict from column G for every row that doesn't
# already carry a change description in column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("changes")

# Clear the AutoFilter criteria on column C ("Labour"). This both drops the
# stored filter criteria from the sheet and un-hides the rows that were
# filtered out. Do this BEFORE touching any cell values so the newly
# revealed rows don't pick up a stray explicit row height when written to.
$null = $ws.Range("A1:K43").AutoFilter(3)

# New header for column H.
$ws.Range("H1").Value = "change type"

# Rows whose G column already evaluates to "No change" but which don't yet
# have a column-H annotation get "No change" filled in.
$noChangeRows = @(4,7,9,11,13,14,18,24,27,28,29,31,32,33,34,36,37,40,41,42)
foreach ($r in $noChangeRows) {
    $ws.Range("H$r").Value = "No change"
}

# Leave the selection where the author left it.
$null = $ws.Range("I5").Select()
